# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '316.95'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '4.44%'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '47.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '10.25%'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.276'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '4.26%'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07928'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '3.27%'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '4.593'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.319'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '31.55%'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.639'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '2.24%'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1279'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '2.48%'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1938'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '4.39%'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.09362'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '3.34%'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.04642'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '11.51%'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1046'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.04%'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001320'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '3.03%'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.04166'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.15%'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005868'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.94%'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.329'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '0.02%'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '3.16%'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.3488'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '4.09%'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.100'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-4.01%'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1395'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-0.26%'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.001321'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '3.25%'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.004196'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-6.56%'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0001353'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '0.45%'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0003548'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-95.23%'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02649'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '8.22%'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05749'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '8.89%'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '83.56%'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.008021'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '4.80%'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '6.54%'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007700'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '4.96%'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.008490'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '13.85%'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3161'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '4.68%'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00006926'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '3.64%'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.36%'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05489'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '34.09%'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.004010'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-4.54%'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002105'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.36%'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002005'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.36%'
